# Apply the "add kvcache per gpu in excel or console output" edit.
#
# Translates the Chinese worksheet/labels to English, adds a new
# "Weight/Single GPU All Layers" column (Q) to the per-operator table,
# zeroes out the dispatch/combine "n" column, and appends two new summary
# rows reporting weight / kv-cache memory per GPU.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Sheet (tab) name
# ---------------------------------------------------------------------
$ws.Name = "Performance Analysis"

# ---------------------------------------------------------------------
# 2. Title cell
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Performance Analysis Report: deepseek_v3 (DECODE)"

# ---------------------------------------------------------------------
# 3. Translate header row (row 3)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Operator Name"
$ws.Range("B3").Value = "Type"
$ws.Range("H3").Value = "Input"
$ws.Range("I3").Value = "Output"
$ws.Range("J3").Value = "Weight"
$ws.Range("K3").Value = "Compute(us)"
$ws.Range("L3").Value = "Memory(us)"
$ws.Range("M3").Value = "Transfer(us)"
$ws.Range("N3").Value = "Single Layer Latency(us)"
$ws.Range("O3").Value = "Total Time(ms)"
$ws.Range("P3").Value = "Percent(%)"

# ---------------------------------------------------------------------
# 4. New column Q: header + width + data rows 4-20
# ---------------------------------------------------------------------
$ws.Range("Q3").Value = "Weight/Single GPU All Layers"
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial($xlPasteFormats)

$ws.Columns.Item(17).ColumnWidth = 11.15

$weightPerGpu = @{
    4  = 923467776
    5  = 2302672896
    6  = 15990784
    7  = 15990784
    8  = 7163871232
    9  = 792723456
    10 = 396361728
    11 = 425721856
    12 = 1702887424
    13 = 851443712
    14 = 1702887424
    15 = 851443712
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
}
foreach ($r in 4..20) {
    $ws.Cells.Item($r, 17).Value = $weightPerGpu[$r]
}
$ws.Range("C4").Copy()
$ws.Range("Q4:Q20").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# 5. dispatch/combine rows: "n" column (D) now reports 0 instead of 7168
# ---------------------------------------------------------------------
$ws.Range("D19").Value = 0
$ws.Range("D20").Value = 0

# ---------------------------------------------------------------------
# 6. Translate the summary labels below the table
# ---------------------------------------------------------------------
$ws.Range("A24").Value = "Compute Time (ms)"
$ws.Range("A25").Value = "Memory Time (ms)"
$ws.Range("A26").Value = "Transfer Time (ms)"
$ws.Range("A27").Value = "Total Time (ms)"
$ws.Range("A30").Value = "Performance Bottleneck"
$ws.Range("B30").Value = "combine (Total Time: 18.804 ms)"
$ws.Range("A33").Value = "Throughput TPS"

# ---------------------------------------------------------------------
# 7. New rows 34/35: weight + kv-cache memory per GPU
# ---------------------------------------------------------------------
$ws.Range("A34").Value = "Weight Memory/Single GPU (GB)"
$ws.Range("B34").Value = 15.968
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial($xlPasteFormats)
$ws.Range("B33").Copy()
$ws.Range("B34").PasteSpecial($xlPasteFormats)

$ws.Range("A35").Value = "KV Cache Memory/Single GPU (GB)"
$ws.Range("B35").Value = 0.000033
$ws.Range("A33").Copy()
$ws.Range("A35").PasteSpecial($xlPasteFormats)
$ws.Range("B35").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 8. Merged title cell now spans through the new column Q
# ---------------------------------------------------------------------
$ws.Range("A1:O1").UnMerge()
$ws.Range("A1:Q1").Merge()
